$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(83,1).Value = "'2021-04-05"
$ws.Cells.Item(83,2).Value = "17:38:49"
$ws.Cells.Item(83,3).Value = "'2.3333"
$ws.Cells.Item(83,4).Value = "'3336.0"
$ws.Cells.Item(83,5).Value = "'2924.0"
$ws.Cells.Item(83,8).Value = "InService"

$ws.Cells.Item(84,1).Value = "'2021-04-05"
$ws.Cells.Item(84,2).Value = "17:40:34"
$ws.Cells.Item(84,3).Value = "'2.3729"
$ws.Cells.Item(84,4).Value = "'5391.0"
$ws.Cells.Item(84,5).Value = "'5733.0"
$ws.Cells.Item(84,8).Value = "InService"

$ws.Cells.Item(85,1).Value = "'2021-04-05"
$ws.Cells.Item(85,2).Value = "17:41:41"
$ws.Cells.Item(85,3).Value = "'2.2951"
$ws.Cells.Item(85,4).Value = "'5792.0"
$ws.Cells.Item(85,5).Value = "'5661.0"
$ws.Cells.Item(85,8).Value = "InService"

$ws.Cells.Item(86,1).Value = "'2021-04-05"
$ws.Cells.Item(86,2).Value = "17:42:41"
$ws.Cells.Item(86,3).Value = "'3.5593"
$ws.Cells.Item(86,4).Value = "'14308.0"
$ws.Cells.Item(86,5).Value = "'12436.0"
$ws.Cells.Item(86,8).Value = "InService"

$ws.Cells.Item(87,1).Value = "'2021-04-05"
$ws.Cells.Item(87,2).Value = "17:43:41"
$ws.Cells.Item(87,3).Value = "'2.623"
$ws.Cells.Item(87,4).Value = "'6926.0"
$ws.Cells.Item(87,5).Value = "'5312.0"
$ws.Cells.Item(87,8).Value = "InService"

$ws.Cells.Item(88,1).Value = "'2021-04-05"
$ws.Cells.Item(88,2).Value = "17:44:41"
$ws.Cells.Item(88,3).Value = "'2.8333"
$ws.Cells.Item(88,4).Value = "'5812.0"
$ws.Cells.Item(88,5).Value = "'4232.0"
$ws.Cells.Item(88,8).Value = "InService"

$ws.Cells.Item(89,1).Value = "'2021-04-05"
$ws.Cells.Item(89,2).Value = "17:45:41"
$ws.Cells.Item(89,3).Value = "'2.8333"
$ws.Cells.Item(89,4).Value = "'7288.0"
$ws.Cells.Item(89,5).Value = "'5856.0"
$ws.Cells.Item(89,8).Value = "InService"

$ws.Cells.Item(90,1).Value = "'2021-04-05"
$ws.Cells.Item(90,2).Value = "17:46:41"
$ws.Cells.Item(90,3).Value = "'2.8333"
$ws.Cells.Item(90,4).Value = "'7810.0"
$ws.Cells.Item(90,5).Value = "'6112.0"
$ws.Cells.Item(90,8).Value = "InService"

$ws.Cells.Item(91,1).Value = "'2021-04-05"
$ws.Cells.Item(91,2).Value = "17:47:41"
$ws.Cells.Item(91,3).Value = "'0.8333"
$ws.Cells.Item(91,4).Value = "'5318.0"
$ws.Cells.Item(91,5).Value = "'3500.0"
$ws.Cells.Item(91,8).Value = "InService"

$ws.Cells.Item(92,1).Value = "'2021-04-05"
$ws.Cells.Item(92,2).Value = "17:48:41"
$ws.Cells.Item(92,3).Value = "'3.0"
$ws.Cells.Item(92,4).Value = "'9102.0"
$ws.Cells.Item(92,5).Value = "'7301.0"
$ws.Cells.Item(92,8).Value = "InService"

$ws.Cells.Item(93,1).Value = "'2021-04-05"
$ws.Cells.Item(93,2).Value = "17:49:41"
$ws.Cells.Item(93,3).Value = "'2.623"
$ws.Cells.Item(93,4).Value = "'5688.0"
$ws.Cells.Item(93,5).Value = "'4540.0"
$ws.Cells.Item(93,8).Value = "InService"

$ws.Cells.Item(94,1).Value = "'2021-04-05"
$ws.Cells.Item(94,2).Value = "17:50:41"
$ws.Cells.Item(94,3).Value = "'1.8644"
$ws.Cells.Item(94,4).Value = "'6976.0"
$ws.Cells.Item(94,5).Value = "'5144.0"
$ws.Cells.Item(94,8).Value = "InService"

$ws.Cells.Item(95,1).Value = "'2021-04-05"
$ws.Cells.Item(95,2).Value = "17:51:41"
$ws.Cells.Item(95,3).Value = "'3.8333"
$ws.Cells.Item(95,4).Value = "'6176.0"
$ws.Cells.Item(95,5).Value = "'4541.0"
$ws.Cells.Item(95,8).Value = "InService"

$ws.Cells.Item(96,1).Value = "'2021-04-05"
$ws.Cells.Item(96,2).Value = "17:52:41"
$ws.Cells.Item(96,3).Value = "'2.6667"
$ws.Cells.Item(96,4).Value = "'5708.0"
$ws.Cells.Item(96,5).Value = "'4332.0"
$ws.Cells.Item(96,8).Value = "InService"

$ws.Cells.Item(97,1).Value = "'2021-04-05"
$ws.Cells.Item(97,2).Value = "17:53:41"
$ws.Cells.Item(97,3).Value = "'2.3333"
$ws.Cells.Item(97,4).Value = "'9732.0"
$ws.Cells.Item(97,5).Value = "'7712.0"
$ws.Cells.Item(97,8).Value = "InService"

$ws.Cells.Item(98,1).Value = "'2021-04-05"
$ws.Cells.Item(98,2).Value = "18:19:24"
$ws.Cells.Item(98,3).Value = "'2.3333"
$ws.Cells.Item(98,4).Value = "'3550.0"
$ws.Cells.Item(98,5).Value = "'3064.0"
$ws.Cells.Item(98,8).Value = "InService"

$ws.Cells.Item(99,1).Value = "'2021-04-05"
$ws.Cells.Item(99,2).Value = "18:19:57"
$ws.Cells.Item(99,3).Value = "'2.2951"
$ws.Cells.Item(99,4).Value = "'5511.0"
$ws.Cells.Item(99,5).Value = "'5870.0"
$ws.Cells.Item(99,8).Value = "InService"

$ws.Cells.Item(100,1).Value = "'2021-04-05"
$ws.Cells.Item(100,2).Value = "18:20:12"
$ws.Cells.Item(100,3).Value = "'2.2951"
$ws.Cells.Item(100,4).Value = "'5511.0"
$ws.Cells.Item(100,5).Value = "'5870.0"
$ws.Cells.Item(100,8).Value = "InService"

